# Fix a few typos on the "Acro" sheet of the acronym/notation workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Acro")

# "Gallium arsenide phosphide" -> "Gallium Arsenide Phosphide"
$ws.Range("B17").Value = "Gallium Arsenide Phosphide"

# "CCR" -> "CCR5"
$ws.Range("A38").Value = "CCR5"

# "Alternating Laser Excitation" -> "Alternating Laser EXcitation"
$ws.Range("B44").Value = "Alternating Laser EXcitation"

# Move the active selection to B45, matching where the author left off editing.
$ws.Activate() | Out-Null
$ws.Range("B45").Select() | Out-Null
